$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A and B, rows 1-32
$values = @(
    @(-0.28018141678364827, 0.27958946117234973),
    @(-0.19029597057718384, 0.18884614509850728),
    @(-0.13912870226378615, 0.13869492901214819),
    @(-0.130694929053913, 0.13002019056455083),
    @(-0.12702019059012937, 0.12470001217433069),
    @(-0.025398229703078812, 0.025099479822760173),
    @(-0.015099479883573075, 0.015037566914913647),
    @(0.00042062892834504595, -0.00043736218022294437),
    @(0.0024373621497075781, -0.0024451257482187572),
    @(0.0044451257178224068, -0.0044448817055879886),
    @(0.0074448816713106325, -0.0074449082654597376),
    @(-0.020872290747329103, 0.020674383703536137),
    @(-0.017174383741944688, 0.017084146742273987),
    @(-0.0090841467988322933, 0.0090545676013649157),
    @(-0.0080545676311629677, 0.0080354576771677344),
    @(-0.0060354577112633478, 0.0060036694823382142),
    @(-0.0040036695170186931, 0.0039999999573980816),
    @(-0.08329349061988367, 0.083118341741212021),
    @(-0.079118341758955157, 0.077793799548817599),
    @(-0.073793799574419339, 0.073396287208479905),
    @(-0.0040058896964687563, 0.0039999999722963864),
    @(-0.045717442857855062, 0.045502414836352756),
    @(-0.040502414864637792, 0.040099736646438622),
    @(-0.020099736735226692, 0.019999999910082167),
    @(-0.097301782500631262, 0.097173908832676403),
    @(-0.094673908865264167, 0.094509464156423917),
    @(-0.092009464190900392, 0.091037197156875749),
    @(-0.089037197196863538, 0.08837146186919842),
    @(-0.081371461934199196, 0.081178305243771121),
    @(-0.021178305517436602, 0.021024757923206572),
    @(-0.014024757992981307, 0.01400143558917577),
    @(-0.0040014356709043852, 0.0039999999414721543)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Swap column widths: column A becomes 15.7109375, column B becomes 16.42578125
$ws.Columns.Item(1).ColumnWidth = 15.7109375
$ws.Columns.Item(2).ColumnWidth = 16.42578125
